$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Cells.Item(24, 7).Value = 1.91
$ws.Cells.Item(24, 8).Value = 3.3
$ws.Cells.Item(24, 9).Value = 4.2
$ws.Cells.Item(24, 10).Value = 2.63
$ws.Cells.Item(24, 12).Value = 4.75
$ws.Cells.Item(24, 13).Value = 1.08
$ws.Cells.Item(24, 14).Value = 8
$ws.Cells.Item(24, 15).Value = 1.4
$ws.Cells.Item(24, 16).Value = 2.75
$ws.Cells.Item(24, 25).Value = 2
$ws.Cells.Item(24, 26).Value = 1.75
$ws.Cells.Item(24, 28).Value = 8.5
$ws.Cells.Item(24, 30).Value = 17
$ws.Cells.Item(24, 35).Value = 17
$ws.Cells.Item(24, 37).Value = 401
$ws.Cells.Item(24, 38).Value = 10
$ws.Cells.Item(24, 39).Value = 21
$ws.Cells.Item(24, 40).Value = 15
$ws.Cells.Item(24, 41).Value = 41
$ws.Cells.Item(24, 43).Value = 41

# Row 25
$ws.Cells.Item(25, 7).Value = 2.3
$ws.Cells.Item(25, 9).Value = 3.4
$ws.Cells.Item(25, 10).Value = 3.2
$ws.Cells.Item(25, 12).Value = 4.33
$ws.Cells.Item(25, 13).Value = 1.11
$ws.Cells.Item(25, 14).Value = 6.5
$ws.Cells.Item(25, 28).Value = 9.5
$ws.Cells.Item(25, 38).Value = 7
$ws.Cells.Item(25, 42).Value = 34

# Row 29
$ws.Cells.Item(29, 15).Value = 1.62
$ws.Cells.Item(29, 16).Value = 2.2

# Row 56
$ws.Cells.Item(56, 7).Value = 1.44
$ws.Cells.Item(56, 8).Value = 3.9
$ws.Cells.Item(56, 19).Value = 2.9
$ws.Cells.Item(56, 20).Value = 1.4
$ws.Cells.Item(56, 25).Value = 2.25
$ws.Cells.Item(56, 26).Value = 1.57
$ws.Cells.Item(56, 34).Value = 8
$ws.Cells.Item(56, 44).Value = 1.54
$ws.Cells.Item(56, 45).Value = 2.43

# Row 70
$ws.Cells.Item(70, 13).Value = 1.06
$ws.Cells.Item(70, 14).Value = 10

# Row 140
$ws.Cells.Item(140, 7).Value = 1.65
$ws.Cells.Item(140, 8).Value = 4.33
$ws.Cells.Item(140, 10).Value = 2.1
$ws.Cells.Item(140, 13).Value = 1.02
$ws.Cells.Item(140, 14).Value = 19
$ws.Cells.Item(140, 17).Value = 1.5
$ws.Cells.Item(140, 18).Value = 2.5
$ws.Cells.Item(140, 19).Value = 1.8
$ws.Cells.Item(140, 20).Value = 2.05
$ws.Cells.Item(140, 21).Value = 2.2
$ws.Cells.Item(140, 22).Value = 1.62
$ws.Cells.Item(140, 23).Value = 1.25
$ws.Cells.Item(140, 24).Value = 3.75
$ws.Cells.Item(140, 25).Value = 1.53
$ws.Cells.Item(140, 26).Value = 2.38
$ws.Cells.Item(140, 27).Value = 11
$ws.Cells.Item(140, 30).Value = 13
$ws.Cells.Item(140, 33).Value = 19
$ws.Cells.Item(140, 36).Value = 34
$ws.Cells.Item(140, 37).Value = 101
$ws.Cells.Item(140, 38).Value = 19
$ws.Cells.Item(140, 39).Value = 29

# Row 144
$ws.Cells.Item(144, 7).Value = 2.1
$ws.Cells.Item(144, 8).Value = 3.15
$ws.Cells.Item(144, 9).Value = 3.35
$ws.Cells.Item(144, 10).Value = 2.62
$ws.Cells.Item(144, 11).Value = 2.1
$ws.Cells.Item(144, 12).Value = 3.8
$ws.Cells.Item(144, 16).Value = 2.9
$ws.Cells.Item(144, 25).Value = 1.72
$ws.Cells.Item(144, 26).Value = 1.9
$ws.Cells.Item(144, 27).Value = 7.3
$ws.Cells.Item(144, 28).Value = 10.25
$ws.Cells.Item(144, 29).Value = 8.5
$ws.Cells.Item(144, 30).Value = 20
$ws.Cells.Item(144, 31).Value = 17
$ws.Cells.Item(144, 33).Value = 9
$ws.Cells.Item(144, 34).Value = 6.1
$ws.Cells.Item(144, 38).Value = 9.75
$ws.Cells.Item(144, 39).Value = 18
$ws.Cells.Item(144, 40).Value = 11.5
$ws.Cells.Item(144, 41).Value = 50
$ws.Cells.Item(144, 42).Value = 32

# Row 150
$ws.Cells.Item(150, 7).Value = 1.95
$ws.Cells.Item(150, 8).Value = 3.3
$ws.Cells.Item(150, 9).Value = 4.2
$ws.Cells.Item(150, 10).Value = 2.6
$ws.Cells.Item(150, 11).Value = 2.1
$ws.Cells.Item(150, 12).Value = 4.5
$ws.Cells.Item(150, 17).Value = 2.05
$ws.Cells.Item(150, 18).Value = 1.8
$ws.Cells.Item(150, 21).Value = 3.5
$ws.Cells.Item(150, 22).Value = 1.29
$ws.Cells.Item(150, 23).Value = 1.44
$ws.Cells.Item(150, 24).Value = 2.63
$ws.Cells.Item(150, 28).Value = 9
$ws.Cells.Item(150, 29).Value = 9
$ws.Cells.Item(150, 30).Value = 17
$ws.Cells.Item(150, 31).Value = 17
$ws.Cells.Item(150, 33).Value = 9
$ws.Cells.Item(150, 37).Value = 251
$ws.Cells.Item(150, 38).Value = 11
$ws.Cells.Item(150, 39).Value = 21
$ws.Cells.Item(150, 40).Value = 15
$ws.Cells.Item(150, 41).Value = 41
$ws.Cells.Item(150, 42).Value = 34

# Row 211
$ws.Cells.Item(211, 7).Value = 2.8
$ws.Cells.Item(211, 9).Value = 2.38
$ws.Cells.Item(211, 10).Value = 3.25
$ws.Cells.Item(211, 11).Value = 2.3
$ws.Cells.Item(211, 12).Value = 3
$ws.Cells.Item(211, 13).Value = 1.03
$ws.Cells.Item(211, 14).Value = 15
$ws.Cells.Item(211, 15).Value = 1.18
$ws.Cells.Item(211, 16).Value = 4.5
$ws.Cells.Item(211, 17).Value = 1.65
$ws.Cells.Item(211, 18).Value = 2.2
$ws.Cells.Item(211, 21).Value = 2.5
$ws.Cells.Item(211, 22).Value = 1.5
$ws.Cells.Item(211, 30).Value = 29
$ws.Cells.Item(211, 34).Value = 7
$ws.Cells.Item(211, 40).Value = 9.5
$ws.Cells.Item(211, 41).Value = 23
$ws.Cells.Item(211, 42).Value = 17

# Row 212
$ws.Cells.Item(212, 7).Value = 1.8
$ws.Cells.Item(212, 8).Value = 3.7
$ws.Cells.Item(212, 9).Value = 4.33
$ws.Cells.Item(212, 10).Value = 2.5
$ws.Cells.Item(212, 15).Value = 1.36
$ws.Cells.Item(212, 16).Value = 3
$ws.Cells.Item(212, 17).Value = 2.1
$ws.Cells.Item(212, 18).Value = 1.7
$ws.Cells.Item(212, 30).Value = 15
$ws.Cells.Item(212, 33).Value = 9

# Row 214
$ws.Cells.Item(214, 7).Value = 2.7
$ws.Cells.Item(214, 9).Value = 2.5
$ws.Cells.Item(214, 10).Value = 3.2
$ws.Cells.Item(214, 12).Value = 3
$ws.Cells.Item(214, 30).Value = 29

# Row 215
$ws.Cells.Item(215, 7).Value = 2.15
$ws.Cells.Item(215, 9).Value = 3.2
$ws.Cells.Item(215, 11).Value = 2.3
$ws.Cells.Item(215, 12).Value = 3.6
$ws.Cells.Item(215, 25).Value = 1.57
$ws.Cells.Item(215, 26).Value = 2.25
$ws.Cells.Item(215, 27).Value = 10
$ws.Cells.Item(215, 28).Value = 12
$ws.Cells.Item(215, 32).Value = 21
$ws.Cells.Item(215, 37).Value = 126
$ws.Cells.Item(215, 38).Value = 13

# Row 216
$ws.Cells.Item(216, 7).Value = 1.42
$ws.Cells.Item(216, 8).Value = 4.75
$ws.Cells.Item(216, 9).Value = 7
$ws.Cells.Item(216, 11).Value = 2.5
$ws.Cells.Item(216, 12).Value = 6.5
$ws.Cells.Item(216, 15).Value = 1.18
$ws.Cells.Item(216, 16).Value = 4.5
$ws.Cells.Item(216, 17).Value = 1.62
$ws.Cells.Item(216, 18).Value = 2.25
$ws.Cells.Item(216, 21).Value = 2.5
$ws.Cells.Item(216, 22).Value = 1.5
$ws.Cells.Item(216, 27).Value = 8
$ws.Cells.Item(216, 31).Value = 11
$ws.Cells.Item(216, 33).Value = 15
$ws.Cells.Item(216, 34).Value = 9
$ws.Cells.Item(216, 38).Value = 19

# Row 217
$ws.Cells.Item(217, 7).Value = 1.91
$ws.Cells.Item(217, 8).Value = 3.75
$ws.Cells.Item(217, 9).Value = 3.75
$ws.Cells.Item(217, 10).Value = 2.5
$ws.Cells.Item(217, 11).Value = 2.38
$ws.Cells.Item(217, 12).Value = 4
$ws.Cells.Item(217, 13).Value = 1.03
$ws.Cells.Item(217, 14).Value = 15
$ws.Cells.Item(217, 15).Value = 1.18
$ws.Cells.Item(217, 16).Value = 4.5
$ws.Cells.Item(217, 17).Value = 1.62
$ws.Cells.Item(217, 18).Value = 2.25
$ws.Cells.Item(217, 21).Value = 2.5
$ws.Cells.Item(217, 22).Value = 1.5
$ws.Cells.Item(217, 23).Value = 1.3
$ws.Cells.Item(217, 24).Value = 3.4
$ws.Cells.Item(217, 27).Value = 10
$ws.Cells.Item(217, 29).Value = 8.5
$ws.Cells.Item(217, 30).Value = 17
$ws.Cells.Item(217, 31).Value = 13
$ws.Cells.Item(217, 32).Value = 21
$ws.Cells.Item(217, 33).Value = 15
$ws.Cells.Item(217, 34).Value = 7
$ws.Cells.Item(217, 36).Value = 34
$ws.Cells.Item(217, 38).Value = 15
$ws.Cells.Item(217, 39).Value = 21
$ws.Cells.Item(217, 40).Value = 13
$ws.Cells.Item(217, 41).Value = 41

# Row 218
$ws.Cells.Item(218, 7).Value = 2.1
$ws.Cells.Item(218, 8).Value = 3.5
$ws.Cells.Item(218, 9).Value = 3.5
$ws.Cells.Item(218, 10).Value = 2.63
$ws.Cells.Item(218, 11).Value = 2.3
$ws.Cells.Item(218, 12).Value = 3.75
$ws.Cells.Item(218, 15).Value = 1.2
$ws.Cells.Item(218, 16).Value = 4.33
$ws.Cells.Item(218, 17).Value = 1.7
$ws.Cells.Item(218, 18).Value = 2.1
$ws.Cells.Item(218, 21).Value = 2.63
$ws.Cells.Item(218, 22).Value = 1.44
$ws.Cells.Item(218, 23).Value = 1.33
$ws.Cells.Item(218, 24).Value = 3.25
$ws.Cells.Item(218, 27).Value = 9.5
$ws.Cells.Item(218, 28).Value = 11
$ws.Cells.Item(218, 30).Value = 19
$ws.Cells.Item(218, 31).Value = 15
$ws.Cells.Item(218, 32).Value = 21
$ws.Cells.Item(218, 33).Value = 13
$ws.Cells.Item(218, 38).Value = 13
$ws.Cells.Item(218, 39).Value = 19
$ws.Cells.Item(218, 41).Value = 41
$ws.Cells.Item(218, 42).Value = 26

# Row 222
$ws.Cells.Item(222, 7).Value = 2.1
$ws.Cells.Item(222, 8).Value = 2.9
$ws.Cells.Item(222, 9).Value = 3.45
$ws.Cells.Item(222, 10).Value = 2.77
$ws.Cells.Item(222, 12).Value = 4
$ws.Cells.Item(222, 15).Value = 1.36
$ws.Cells.Item(222, 16).Value = 2.67
$ws.Cells.Item(222, 17).Value = 2.12
$ws.Cells.Item(222, 21).Value = 3.7
$ws.Cells.Item(222, 22).Value = 1.19
$ws.Cells.Item(222, 23).Value = 1.42
$ws.Cells.Item(222, 24).Value = 2.45
$ws.Cells.Item(222, 25).Value = 1.82
$ws.Cells.Item(222, 27).Value = 6.4
$ws.Cells.Item(222, 28).Value = 9.75
$ws.Cells.Item(222, 29).Value = 8.75
$ws.Cells.Item(222, 30).Value = 21
$ws.Cells.Item(222, 31).Value = 19
$ws.Cells.Item(222, 33).Value = 7.5
$ws.Cells.Item(222, 34).Value = 5.8
$ws.Cells.Item(222, 35).Value = 15
$ws.Cells.Item(222, 37).Value = 700
$ws.Cells.Item(222, 38).Value = 9.25
$ws.Cells.Item(222, 39).Value = 18.5
$ws.Cells.Item(222, 40).Value = 12
$ws.Cells.Item(222, 41).Value = 55
$ws.Cells.Item(222, 42).Value = 35
